$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (datetime) corrections for existing rows ---
$ws.Range("A4").Value  = 35156.45833333334
$ws.Range("A16").Value = 35521.45833333334
$ws.Range("A28").Value = 35886.45833333334
$ws.Range("A40").Value = 36251.45833333334
$ws.Range("A52").Value = 36617.45833333334
$ws.Range("A57").Value = 36770.41666666666
$ws.Range("A76").Value = 37347.45833333334
$ws.Range("A88").Value = 37712.45833333334
$ws.Range("A100").Value = 38078.45833333334
$ws.Range("A112").Value = 38443.45833333334
$ws.Range("A130").Value = 38991.45833333334

# --- Updated open/high/low/close values for existing rows ---
$ws.Range("C318:F318").Value = 6614488000000
$ws.Range("C319:F319").Value = 6708871000000
$ws.Range("C320:F320").Value = 6808405000000
$ws.Range("C327:F327").Value = 7965913000000
$ws.Range("C329:F329").Value = 8140532000000

# --- New row 330, append new data point ---
$ws.Range("A329:G329").Copy()
$ws.Range("A330:G330").PasteSpecial(-4122)  # xlPasteFormats, keep same styling as row above

$ws.Range("A330").Value = 45078.41666666666
$ws.Range("B330").Value = "ECONOMICS:EGM2"
$ws.Range("C330:F330").Value = 8248190000000
$ws.Range("G330").Value = 0
